# Week3.pptx edit: split the "Variable names" bullet on slide 7 into three
# runs so the wording reads "Must start with a letter or an underscore,
# followed by any sequence of letters, digits, and underscores (invalid
# names: 2days, five/three)".

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(7)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$run1 = "Must start with a letter or "
$run2 = "an underscore, followed "
$run3 = "by any sequence of letters, digits, and underscores (invalid names: 2days, five/three)"
$newText = $run1 + $run2 + $run3

# Locate paragraph 9 ("Must start with a letter, follow by any sequence...")
# inside the "Variable names" content placeholder and replace its text.
$para = $tr.Paragraphs(9, 1)
$paraStart = $para.Start

$whole = $tr.Characters($paraStart, $para.Length)
$whole.Text = $newText

# Re-assign the middle segment on its own so PowerPoint keeps it as its own
# run (matching the three-run split introduced by the edit) instead of
# merging it back with its neighbours.
$middle = $tr.Characters($paraStart + $run1.Length, $run2.Length)
$middle.Text = $run2
